# The source data pipeline now emits three additional general-college-subject
# columns (history, electives, cs) ahead of the existing "arts" column, and
# normalizes the "Unknown" placeholder strings to lowercase "unknown".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank columns at R:T, pushing the existing
# "general_college_subjects.arts" column (and everything after it) right by three.
$ws.Range("R:T").Insert()

# Header row (row 1) labels for the newly inserted columns.
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# Data row (row 2) values for the newly inserted columns.
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0

# Normalize the "Unknown" placeholders to lowercase "unknown".
$ws.Range("D2:J2").Value = "unknown"
